{"js": "// Add \" a Tropy\" right after the existing \"Jazykov\u00e9 prost\u0159edky\" heading text,\n// so the heading reads \"Jazykov\u00e9 prost\u0159edky a Tropy\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\nconst target = paragraphs.items.find(p => p.text.trim() === \"Jazykov\u00e9 prost\u0159edky\");\n\nif (target) {\n  // Insert the new text immediately at the end of the paragraph (after the\n  // existing \"Jazykov\u00e9 prost\u0159edky\" run), preserving the leading space so the\n  // final text reads \"Jazykov\u00e9 prost\u0159edky a Tropy\".\n  const endRange = target.getRange(\"End\");\n  endRange.insertText(\" a Tropy\", Word.InsertLocation.end);\n  await context.sync();\n}\n", "ps1": "# Add \" a Tropy\" right after the existing \"Jazykov\u00e9 prost\u0159edky\" heading text,\n# so the heading reads \"Jazykov\u00e9 prost\u0159edky a Tropy\".\n\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $r = $p.Range\n    $text = $r.Text.TrimEnd([char]13, [char]7).Trim()\n    if ($text -eq \"Jazykov\u00e9 prost\u0159edky\") {\n        $endRange = $p.Range\n        $endRange.Collapse(0) | Out-Null         # wdCollapseEnd\n        $endRange.MoveEnd(1, -1) | Out-Null      # move before the paragraph mark\n        $endRange.InsertAfter(\" a Tropy\") | Out-Null\n        break\n    }\n}\n"}
